# Apply uniform paragraph spacing (After: 6pt, Line spacing: 1.5 lines)
# to every paragraph in the document body, matching:
#   <w:spacing w:after="120" w:line="360" w:lineRule="auto"/>
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $pf = $p.Range.ParagraphFormat
    $pf.LineSpacingRule = 1   # wdLineSpace1pt5 -> w:line="360" w:lineRule="auto"
    $pf.SpaceAfter = 6        # 6pt -> w:after="120"
}

Write-Output "Applied spacing (after=6pt, line=1.5) to $($d.Paragraphs.Count) paragraphs"
